$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Vehículos en el hogar"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Municipio código"
$ws.Range("G1").Value = "Municipio nombre"

# --- Row 2 ---
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:vehiculos-en-el-hogar"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# --- Row 3 ---
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"

# --- Row 4 ---
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Municipio"
